# Applies the changes described by the diff between before.xlsx and the
# author's updated CasosColombia.xlsx:
#   - A handful of scattered cells toggle between a numeric value and the
#     literal text "NaN" (shared string already present in the workbook).
#   - A brand-new data row (173, date 2020-08-15 / serial 44067) is appended
#     at the bottom of the sheet, continuing the same daily series as row 172.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Scattered cell edits: numeric <-> the text "NaN"
# ---------------------------------------------------------------------

# Cells that become the text "NaN" (were plain numbers before)
$toNaN = @("I9","I10","L25","L26","L27","CM31","BP34","CM80","CM81","CH88","AP89","H91","AP93","BQ107","BQ108","J115","AF127")
foreach ($addr in $toNaN) {
    $ws.Range($addr).Value = "NaN"
}

# Cells that become real numbers (were the text "NaN" before)
$ws.Range("BP24").Value = 1
$ws.Range("AK58").Value = 6
$ws.Range("AK59").Value = 6
$ws.Range("AK60").Value = 6
$ws.Range("AK61").Value = 6
$ws.Range("AK62").Value = 6
$ws.Range("AK63").Value = 6
$ws.Range("AK64").Value = 6
$ws.Range("AK65").Value = 6
$ws.Range("AK66").Value = 6
$ws.Range("AK67").Value = 6
$ws.Range("AK68").Value = 6
$ws.Range("AK69").Value = 6
$ws.Range("AK70").Value = 6
$ws.Range("AK71").Value = 6
$ws.Range("AK72").Value = 6
$ws.Range("AK98").Value = 74
$ws.Range("AK99").Value = 77

# ---------------------------------------------------------------------
# 2. Append new row 173 (same column layout/format as row 172)
# ---------------------------------------------------------------------

# Copy row 172 down into row 173 first so the new row inherits the exact
# same per-column styling (date style in A, the s="4" block from BS..DQ,
# plain cells elsewhere) as the rest of the table.
$ws.Range("A172:DX172").Copy()
$ws.Range("A173:DX173").PasteSpecial(-4104)
$excel.CutCopyMode = $false

# Now overwrite row 173 with the real values for 2020-08-15.
$row173 = @{
    "A173" = 44067;   "B173" = 551696; "C173" = 2690;   "D173" = 71563;
    "E173" = 63069;   "F173" = 192654; "G173" = 23694;  "H173" = 3001;
    "I173" = 2386;    "J173" = 5069;   "K173" = 4329;   "L173" = 8067;
    "M173" = 3664;    "N173" = 17498;  "O173" = 19626;  "P173" = 4471;
    "Q173" = 3443;    "R173" = 11444;  "S173" = 6089;   "T173" = 12846;
    "U173" = 9134;    "V173" = 2476;   "W173" = 937;    "X173" = 4791;
    "Y173" = 14186;   "Z173" = 10505;  "AA173" = 5691;  "AB173" = 43598;
    "AC173" = 888;    "AD173" = 126;   "AE173" = 216;   "AF173" = 440;
    "AG173" = 49;     "AH173" = 29;    "AI173" = 236;   "AJ173" = 1940;
    "AK173" = 2387;   "AL173" = 35532; "AM173" = 6031;  "AN173" = 2381;
    "AO173" = 34340;  "AP173" = 843;   "AQ173" = 19415; "AR173" = 1415;
    "AS173" = 6370;   "AT173" = 1402;  "AU173" = 1543;  "AV173" = 3277;
    "AW173" = 1450;   "AX173" = 925;   "AY173" = 2458;  "AZ173" = 2573;
    "BA173" = 41035;  "BB173" = 11098; "BC173" = 1953;  "BD173" = 6766;
    "BE173" = 3103;   "BF173" = 274;   "BG173" = 1375;  "BH173" = 2549;
    "BI173" = 727;    "BJ173" = 1950;  "BK173" = 7747;  "BL173" = 7584;
    "BM173" = 7349;   "BN173" = 13618; "BO173" = 1846;  "BP173" = 769;
    "BQ173" = 5530;   "BR173" = 5077;  "BS173" = 5903;  "BT173" = 1265;
    "BU173" = 1351;   "BV173" = 2336;  "BW173" = 2623;  "BX173" = 708;
    "BY173" = 3890;   "BZ173" = 2227;  "CA173" = 1083;  "CB173" = 619;
    "CC173" = 1816;   "CD173" = 1759;  "CE173" = 987;   "CF173" = 829;
    "CG173" = 4224;   "CH173" = 1117;  "CI173" = 1064;  "CJ173" = 1103;
    "CK173" = 1392;   "CL173" = 1270;  "CM173" = 1312;  "CN173" = 1045;
    "CO173" = 983;    "CP173" = 1052;  "CQ173" = 528;   "CR173" = 2855;
    "CS173" = 870;    "CT173" = 743;   "CU173" = 674;   "CV173" = 1136;
    "CW173" = 1031;   "CX173" = 574;   "CY173" = 686;   "CZ173" = 773;
    "DA173" = 1075;   "DB173" = 878;   "DC173" = 1002;  "DD173" = 782;
    "DE173" = 311;    "DF173" = 324;   "DG173" = 629;   "DH173" = 541;
    "DI173" = 389;    "DJ173" = 528;   "DK173" = 308;   "DL173" = 536;
    "DM173" = 693;    "DN173" = 502;   "DO173" = 471;   "DP173" = 345;
    "DQ173" = 508;    "DR173" = 110745; "DS173" = 232225; "DT173" = 8387;
    "DU173" = 99556;  "DV173" = 64437; "DW173" = 23369; "DX173" = 7594;
}

foreach ($addr in $row173.Keys) {
    $ws.Range($addr).Value = $row173[$addr]
}

# ---------------------------------------------------------------------
# 3. Restore the view's active cell to the new last cell of the sheet,
#    matching the author's final selection before saving.
# ---------------------------------------------------------------------
$ws.Range("DX173").Select()
